$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source Data")

# The first iteration block's last row (Week 7 / Iteration 4, row 6) is being
# moved into the second iteration block as its first row (row 9), replacing
# the old "StuVac" row there. Row 6 is then cleared out (values + formatting)
# leaving only the release-number cell (E6) behind.
$ws.Range("A6:D6").Copy($ws.Range("A9:D9"))
$ws.Range("A6:D6").Clear()

# The pivot cache's source range grew from B1:D6 to B1:D9 on "Source Data".
$pt = $ws.Parent.Worksheets.Item("By Product").PivotTables().Item(1)
$pt.SourceData = "Source Data!B1:D9"

# Move the active cell / selection to D3.
$ws.Activate() | Out-Null
$ws.Range("D3").Select() | Out-Null
